# Add scores from 28 DEC 2024 to the "Score Cards" sheet / Table1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score Cards")
$lo = $ws.ListObjects.Item("Table1")

# Adding a ListRow expands the Excel Table (ref + autoFilter) just like
# typing a new row directly below the table in the UI would.
$newRowObj = $lo.ListRows.Add()
$row = $newRowObj.Index + 1   # +1 to account for the header row
$prevRow = $row - 1

# Copy the formatting from the row above down into the new row so the
# new row matches the existing table styling (date format, alignment...).
$ws.Range("A$prevRow`:W$prevRow").Copy()
$ws.Range("A$row`:W$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = "Ocean View"
$ws.Cells.Item($row, 2).Value = (Get-Date -Year 2024 -Month 12 -Day 28).Date

# Hole-by-hole scores (1-18), followed by GIR, Putts and Fairways.
$scores = @(6, 5, 6, 3, 7, 4, 6, 5, 5, 3, 4, 5, 4, 4, 6, 6, 4, 4, 4, 36, 4)
$col = 3
foreach ($val in $scores) {
    $ws.Cells.Item($row, $col).Value = $val
    $col++
}

$ws.Activate()
[void]$ws.Range("P13").Select()

$wb.Save()
